$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value2 = 171.55556
$ws.Range("I9").Value2 = 182.33333
$ws.Range("J9").Value2 = 150
$ws.Range("K9").Value2 = 182.33333
$ws.Range("L9").Value2 = 150
$ws.Range("M9").Value2 = -13.33332999999999
$ws.Range("N9").Value2 = -488
$ws.Range("H39").Value2 = 872.2857
$ws.Range("I39").Value2 = 785.53845
$ws.Range("J39").Value2 = 2000
$ws.Range("K39").Value2 = 2356.61535
$ws.Range("L39").Value2 = 6000
$ws.Range("M39").Value2 = -2060.61535
$ws.Range("N39").Value2 = -6592
$ws.Range("H40").Value2 = 4381.5
$ws.Range("I40").Value2 = 3223.75
$ws.Range("K40").Value2 = 3223.75
$ws.Range("M40").Value2 = -3048.75
$ws.Range("H53").Value2 = 1288.1333
$ws.Range("J53").Value2 = 1051.8
$ws.Range("L53").Value2 = 1051.8
$ws.Range("N53").Value2 = -2325.8
$ws.Range("H98").Value2 = 1634.0625
$ws.Range("I98").Value2 = 1180.8334
$ws.Range("J98").Value2 = 2993.75
$ws.Range("K98").Value2 = 1180.8334
$ws.Range("L98").Value2 = 2993.75
$ws.Range("M98").Value2 = 317.1666
$ws.Range("N98").Value2 = -5989.75
$ws.Range("H113").Value2 = 7238.9
$ws.Range("I113").Value2 = 7161.25
$ws.Range("K113").Value2 = 7161.25
$ws.Range("M113").Value2 = -3907.25
$ws.Range("H122").Value2 = 1634.0625
$ws.Range("I122").Value2 = 1180.8334
$ws.Range("J122").Value2 = 2993.75
$ws.Range("K122").Value2 = 3542.5002
$ws.Range("L122").Value2 = 8981.25
$ws.Range("M122").Value2 = -1092.5002
$ws.Range("N122").Value2 = -13881.25
$ws.Range("H132").Value2 = 45152.824
$ws.Range("I132").Value2 = 47127.953
$ws.Range("K132").Value2 = 141383.859
$ws.Range("M132").Value2 = -138853.859

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 4231.864
$ws.Range("I2").Value2 = 4734.871
$ws.Range("J2").Value2 = 3032.3845
$ws.Range("K2").Value2 = 4734.871
$ws.Range("L2").Value2 = 3032.3845
$ws.Range("M2").Value2 = -4621.871
$ws.Range("N2").Value2 = -3258.3845
$ws.Range("H32").Value2 = 1262.6428
$ws.Range("I32").Value2 = 1262.6428
$ws.Range("K32").Value2 = 1262.6428
$ws.Range("M32").Value2 = -975.6428000000001
$ws.Range("H95").Value2 = 60000
$ws.Range("J95").Value2 = 60000
$ws.Range("L95").Value2 = 60000
$ws.Range("N95").Value2 = -65492
$ws.Range("H116").Value2 = 4231.864
$ws.Range("I116").Value2 = 4734.871
$ws.Range("J116").Value2 = 3032.3845
$ws.Range("K116").Value2 = 4734.871
$ws.Range("L116").Value2 = 3032.3845
$ws.Range("M116").Value2 = -2440.871
$ws.Range("N116").Value2 = -7620.3845
$ws.Range("H132").Value2 = 20837416
$ws.Range("I132").Value2 = 2358.3157
$ws.Range("K132").Value2 = 7074.9471
$ws.Range("M132").Value2 = -4544.9471

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 4231.864
$ws.Range("I3").Value2 = 4734.871
$ws.Range("J3").Value2 = 3032.3845
$ws.Range("K3").Value2 = 4734.871
$ws.Range("L3").Value2 = 3032.3845
$ws.Range("M3").Value2 = -4620.871
$ws.Range("N3").Value2 = -3260.3845
$ws.Range("H44").Value2 = 99999
$ws.Range("J44").Value2 = 99999
$ws.Range("L44").Value2 = 99999
$ws.Range("N44").Value2 = -100993
$ws.Range("H105").Value2 = 2466.7896
$ws.Range("I105").Value2 = 2739.5
$ws.Range("J105").Value2 = 2163.7778
$ws.Range("K105").Value2 = 2739.5
$ws.Range("L105").Value2 = 2163.7778
$ws.Range("M105").Value2 = -992.5
$ws.Range("N105").Value2 = -5657.7778
$ws.Range("H123").Value2 = 49999
$ws.Range("J123").Value2 = 49999
$ws.Range("L123").Value2 = 49999
$ws.Range("N123").Value2 = -59799
$ws.Range("H134").Value2 = 15628386
$ws.Range("I134").Value2 = 19233942
$ws.Range("K134").Value2 = 57701826
$ws.Range("M134").Value2 = -57699291

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value2 = 2502149.2
$ws.Range("I10").Value2 = 2502149.2
$ws.Range("K10").Value2 = 2502149.2
$ws.Range("M10").Value2 = -2502010.2
$ws.Range("H31").Value2 = 2430.8635
$ws.Range("J31").Value2 = 4176.5
$ws.Range("L31").Value2 = 4176.5
$ws.Range("N31").Value2 = -4766.5
$ws.Range("H34").Value2 = 2430.8635
$ws.Range("J34").Value2 = 4176.5
$ws.Range("L34").Value2 = 4176.5
$ws.Range("N34").Value2 = -4580.5
$ws.Range("H56").Value2 = 30103
$ws.Range("J56").Value2 = 30103
$ws.Range("L56").Value2 = 30103
$ws.Range("N56").Value2 = -31793

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value2 = 82.28846
$ws.Range("I40").Value2 = 58.117645
$ws.Range("K40").Value2 = 232.47058
$ws.Range("M40").Value2 = -163.47058
$ws.Range("H56").Value2 = 8895.333000000001
$ws.Range("I56").Value2 = 8895.333000000001
$ws.Range("K56").Value2 = 8895.333000000001
$ws.Range("M56").Value2 = -8365.333000000001
$ws.Range("H62").Value2 = 14999.5
$ws.Range("I62").Value2 = 0
$ws.Range("J62").Value2 = 14999.5
$ws.Range("K62").Value2 = 0
$ws.Range("L62").Value2 = 44998.5
$ws.Range("N62").Value2 = -46370.5
$ws.Range("H65").Value2 = 14999.5
$ws.Range("I65").Value2 = 0
$ws.Range("J65").Value2 = 14999.5
$ws.Range("K65").Value2 = 0
$ws.Range("L65").Value2 = 134995.5
$ws.Range("N65").Value2 = -141859.5
$ws.Range("H69").Value2 = 3243.077
$ws.Range("I69").Value2 = 2200
$ws.Range("J69").Value2 = 3706.6667
$ws.Range("K69").Value2 = 6600
$ws.Range("L69").Value2 = 11120.0001
$ws.Range("M69").Value2 = -5789
$ws.Range("N69").Value2 = -12742.0001
$ws.Range("H72").Value2 = 3243.077
$ws.Range("I72").Value2 = 2200
$ws.Range("J72").Value2 = 3706.6667
$ws.Range("K72").Value2 = 19800
$ws.Range("L72").Value2 = 33360.0003
$ws.Range("M72").Value2 = -15744
$ws.Range("N72").Value2 = -41472.0003
$ws.Range("H122").Value2 = 51955.1
$ws.Range("I122").Value2 = 999
$ws.Range("K122").Value2 = 8991
$ws.Range("M122").Value2 = -6541
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value2 = 2651.4546
$ws.Range("I126").Value2 = 0
$ws.Range("J126").Value2 = 2651.4546
$ws.Range("K126").Value2 = 0
$ws.Range("L126").Value2 = 7954.3638
$ws.Range("N126").Value2 = -12894.3638
$ws.Range("H132").Value2 = 6562
$ws.Range("I132").Value2 = 6249
$ws.Range("K132").Value2 = 18747
$ws.Range("M132").Value2 = -16217
$ws.Range("M126").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value2 = 5450
$ws.Range("I45").Value2 = 5450
$ws.Range("K45").Value2 = 5450
$ws.Range("M45").Value2 = -5043
$ws.Range("H98").Value2 = 69677.5
$ws.Range("J98").Value2 = 69677.5
$ws.Range("L98").Value2 = 69677.5
$ws.Range("N98").Value2 = -75667.5
$ws.Range("H132").Value2 = 3186.3333
$ws.Range("I132").Value2 = 3247.4285
$ws.Range("K132").Value2 = 9742.2855
$ws.Range("M132").Value2 = -7212.2855

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value2 = 1761.2
$ws.Range("J96").Value2 = 1804.7142
$ws.Range("L96").Value2 = 1804.7142
$ws.Range("N96").Value2 = -4550.7142
$ws.Range("H122").Value2 = 2989.7646
$ws.Range("I122").Value2 = 3243.8333
$ws.Range("J122").Value2 = 2380
$ws.Range("K122").Value2 = 9731.499899999999
$ws.Range("L122").Value2 = 7140
$ws.Range("M122").Value2 = -7281.499899999999
$ws.Range("N122").Value2 = -12040
